# Weekly update: insert two new price records (2021-09-10) for "Poroto verde"
# at Vega Central Mapocho de Santiago, right before the existing block of
# records that currently starts at row 214. This pushes the dimension from
# A1:R239 to A1:R241.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 214-215, shifting old rows 214:239 down to 216:241
$ws.Range('A214:A215').EntireRow.Insert()

# New row 214
$ws.Range('A214').Value = 9
$ws.Range('B214').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C214').Value = 'Metropolitana'
$ws.Range('D214').Value = 44449
$ws.Range('E214').Value = 13
$ws.Range('F214').Value = 100112031
$ws.Range('G214').Value = 'Poroto verde'
$ws.Range('H214').Value = 'Magnum'
$ws.Range('I214').Value = 'Primera'
$ws.Range('J214').Value = 25
$ws.Range('K214').Value = 37000
$ws.Range('L214').Value = 38000
$ws.Range('M214').Value = 37480
$ws.Range('N214').Value = '$/malla 25 kilos'
$ws.Range('O214').Value = 'Perú'
$ws.Range('P214').Value = 1499
$ws.Range('Q214').Value = 25
$ws.Range('R214').Value = 'Hortaliza'

# New row 215
$ws.Range('A215').Value = 9
$ws.Range('B215').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C215').Value = 'Metropolitana'
$ws.Range('D215').Value = 44449
$ws.Range('E215').Value = 13
$ws.Range('F215').Value = 100112031
$ws.Range('G215').Value = 'Poroto verde'
$ws.Range('H215').Value = 'Sin especificar'
$ws.Range('I215').Value = 'Primera'
$ws.Range('J215').Value = 18
$ws.Range('K215').Value = 36000
$ws.Range('L215').Value = 37000
$ws.Range('M215').Value = 36500
$ws.Range('N215').Value = '$/malla 25 kilos'
$ws.Range('O215').Value = 'Perú'
$ws.Range('P215').Value = 1460
$ws.Range('Q215').Value = 25
$ws.Range('R215').Value = 'Hortaliza'
